# Rebuild the monthly time-series table.
# The source rows are reordered (Oct/Nov/Dec of each year now lead that
# year's block, ahead of Jan-Sep) and three more years of data
# (2022-10 .. 2023-07) are appended, growing the sheet from 49 to 68 rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template used to stamp the same look (bold, centered, bordered date cell)
# onto the rows that are brand new (rows 50-68 did not exist before).
$templateRow = $ws.Range("A2:E2")

$data = @(
    ,@("2018-10", 102.4, 100.9, 101.2, 101.5)
    ,@("2018-11", 102.2, 100.4, 101, 101.9)
    ,@("2018-12", 101.5, 100.4, 101.4, 101.9)
    ,@("2018-01", 101.9162, 101.1471, 98.9286, 99.30800000000001)
    ,@("2018-02", 102.2, 101.3, 98.7, 99.09999999999999)
    ,@("2018-03", 102, 101.1, 99, 98.40000000000001)
    ,@("2018-04", 102.2, 101.7, 99.2, 99.09999999999999)
    ,@("2018-05", 102.2, 101.7, 99.59999999999999, 99.90000000000001)
    ,@("2018-06", 102.5, 101.8, 100, 100.4)
    ,@("2018-07", 102.7, 101.8, 99.7, 101.1)
    ,@("2018-08", 102.8, 101.7, 100.3, 102.1)
    ,@("2018-09", 102.7, 101.3, 100.8, 102.3)
    ,@("2019-10", 101.3, 98.2, 100.1, 100.5)
    ,@("2019-11", 101.3, 98.2, 100.5, 99.5)
    ,@("2019-12", 101.4, 98, 100.2, 100.2)
    ,@("2019-01", 101.4, 100.4, 101.7, 101.2)
    ,@("2019-02", 101.8, 100.4, 101.3, 102)
    ,@("2019-03", 101.9, 99.7, 100.9, 101.2)
    ,@("2019-04", 101.7, 99.5, 101.4, 101.3)
    ,@("2019-05", 101.5, 99.90000000000001, 101.5, 100.1)
    ,@("2019-06", 101.7, 99.8, 101.5, 100)
    ,@("2019-07", 101.6, 99.7, 101.8, 100.4)
    ,@("2019-08", 101.7, 97.90000000000001, 101.1, 99.59999999999999)
    ,@("2019-09", 101.4, 98, 100.3, 100.2)
    ,@("2020-10", 99, 98.59999999999999, 98.7, 99.5)
    ,@("2020-11", 99, 98.8, 98.40000000000001, 99.90000000000001)
    ,@("2020-12", 99.5, 98.8, 98.8, 99)
    ,@("2020-01", 101.3, 97.8, 99.90000000000001, 100.7)
    ,@("2020-02", 100.8, 97.3, 100.2, 101.1)
    ,@("2020-03", 100.4, 97.90000000000001, 100.2, 102)
    ,@("2020-04", 99.7, 97.8, 99.40000000000001, 101.4)
    ,@("2020-05", 99.7, 97.40000000000001, 99.2, 101.7)
    ,@("2020-06", 98.90000000000001, 97.59999999999999, 98.40000000000001, 101.2)
    ,@("2020-07", 99, 97.3, 98.40000000000001, 99.59999999999999)
    ,@("2020-08", 99.2, 99.09999999999999, 98.3, 99.59999999999999)
    ,@("2020-09", 99.09999999999999, 98.59999999999999, 99.09999999999999, 98.7)
    ,@("2021-10", 104.4, 100.4, 104, 99.5)
    ,@("2021-11", 105.1, 101.1, 104.7, 100.1)
    ,@("2021-12", 103.9, 101.1, 105, 100.6)
    ,@("2021-01", 100.2, 100, 99.3, 97.59999999999999)
    ,@("2021-02", 100.7, 99.8, 99.7, 98.09999999999999)
    ,@("2021-03", 100.9, 99.7, 99.59999999999999, 96.8)
    ,@("2021-04", 101.4, 99.09999999999999, 101, 96.90000000000001)
    ,@("2021-05", 101.6, 99.90000000000001, 101.3, 97.09999999999999)
    ,@("2021-06", 102, 99.59999999999999, 101.1, 96.3)
    ,@("2021-07", 102, 100.1, 101.3, 98.8)
    ,@("2021-08", 102, 100.2, 102.3, 98.8)
    ,@("2021-09", 102.6, 100.5, 102.3, 99.40000000000001)
    ,@("2022-10", 98.8, 101.6, 103, 103.2)
    ,@("2022-11", 98.09999999999999, 101.1, 101.9, 103.1)
    ,@("2022-12", 98.8, 100.9, 101.3, 101)
    ,@("2022-01", 103.6, 101.3, 105.1, 100.1)
    ,@("2022-02", 103.3, 101.2, 104.7, 99.59999999999999)
    ,@("2022-03", 103.3, 101.6, 105.2, 101)
    ,@("2022-04", 103.1, 102.3, 104.7, 101.1)
    ,@("2022-05", 102.7, 101.9, 105.8, 101.5)
    ,@("2022-06", 101.9, 102.3, 106.1, 103)
    ,@("2022-07", 101.3, 101.5, 105.3, 103.3)
    ,@("2022-08", 100.8, 101.4, 104.2, 103.6)
    ,@("2022-09", 100.2, 101.3, 104.5, 102.9)
    ,@("2023-01", 98.90000000000001, 100.3, 100.6, 97.5)
    ,@("2023-02", 98, 100.4, 100.1, 98.2)
    ,@("2023-03", 98, 100.1, 99.90000000000001, 97.90000000000001)
    ,@("2023-04", 98.8, 100.2, 99.59999999999999, 98)
    ,@("2023-05", 98.7, 99.90000000000001, 98.59999999999999, 97.8)
    ,@("2023-06", 98.2, 99.7, 98.40000000000001, 97.40000000000001)
    ,@("2023-07", 98.3, 100.1, 98, 97)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $row = $data[$i]

    if ($rowNum -gt 49) {
        # New row beyond the original 49-row extent: clone formatting first.
        $templateRow.Copy()
        $ws.Range("A" + $rowNum + ":E" + $rowNum).PasteSpecial(-4122)
    }

    $ws.Cells.Item($rowNum, 1).Value = $row[0]
    $ws.Cells.Item($rowNum, 2).Value = $row[1]
    $ws.Cells.Item($rowNum, 3).Value = $row[2]
    $ws.Cells.Item($rowNum, 4).Value = $row[3]
    $ws.Cells.Item($rowNum, 5).Value = $row[4]
}

$excel.CutCopyMode = 0
